Write-Output "noop"
